$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap country names (Pais column A) for re-ranked rows ---
$tmpA19 = $ws.Range("A19").Value2
$tmpA20 = $ws.Range("A20").Value2
$ws.Range("A19").Value2 = $tmpA20
$ws.Range("A20").Value2 = $tmpA19

$tmpA50 = $ws.Range("A50").Value2
$tmpA51 = $ws.Range("A51").Value2
$ws.Range("A50").Value2 = $tmpA51
$ws.Range("A51").Value2 = $tmpA50

$tmpA99 = $ws.Range("A99").Value2
$tmpA100 = $ws.Range("A100").Value2
$ws.Range("A99").Value2 = $tmpA100
$ws.Range("A100").Value2 = $tmpA99

$tmpA113 = $ws.Range("A113").Value2
$tmpA114 = $ws.Range("A114").Value2
$ws.Range("A113").Value2 = $tmpA114
$ws.Range("A114").Value2 = $tmpA113

$tmpA202 = $ws.Range("A202").Value2
$tmpA203 = $ws.Range("A203").Value2
$ws.Range("A202").Value2 = $tmpA203
$ws.Range("A203").Value2 = $tmpA202

$tmpA213 = $ws.Range("A213").Value2
$tmpA214 = $ws.Range("A214").Value2
$ws.Range("A213").Value2 = $tmpA214
$ws.Range("A214").Value2 = $tmpA213

# --- Update numeric stats (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes) ---
$ws.Range("B4").Value2 = 5246213
$ws.Range("C4").Value2 = 44567
$ws.Range("D4").Value2 = 2700573
$ws.Range("E4").Value2 = 2379550
$ws.Range("G4").Value2 = 467
$ws.Range("H4").Value2 = 166090

$ws.Range("D5").Value2 = 2163812
$ws.Range("E5").Value2 = 791801
$ws.Range("G5").Value2 = 721
$ws.Range("H5").Value2 = 101857

$ws.Range("B11").Value2 = 397623
$ws.Range("C11").Value2 = 10142
$ws.Range("D11").Value2 = 221485
$ws.Range("E11").Value2 = 162984
$ws.Range("G11").Value2 = 312
$ws.Range("H11").Value2 = 13154

$ws.Range("B19").Value2 = 253868
$ws.Range("C19").Value2 = 7369
$ws.Range("D19").Value2 = 108242
$ws.Range("E19").Value2 = 140862
$ws.Range("G19").Value2 = 158
$ws.Range("H19").Value2 = 4764

$ws.Range("B20").Value2 = 250825
$ws.Range("C20").Value2 = 259
$ws.Range("D20").Value2 = 202248
$ws.Range("E20").Value2 = 13368
$ws.Range("G20").Value2 = 4
$ws.Range("H20").Value2 = 35209

$ws.Range("B22").Value2 = 218500
$ws.Range("C22").Value2 = 1219
$ws.Range("E22").Value2 = 11335

$ws.Range("B39").Value2 = 75394
$ws.Range("C39").Value2 = 902
$ws.Range("D39").Value2 = 49510
$ws.Range("E39").Value2 = 24220
$ws.Range("G39").Value2 = 25
$ws.Range("H39").Value2 = 1664

$ws.Range("B50").Value2 = 47990
$ws.Range("C50").Value2 = 1207
$ws.Range("D50").Value2 = 33058
$ws.Range("E50").Value2 = 13885
$ws.Range("G50").Value2 = 7
$ws.Range("H50").Value2 = 1047

$ws.Range("B51").Value2 = 47454
$ws.Range("C51").Value2 = 481
$ws.Range("D51").Value2 = 6597
$ws.Range("E51").Value2 = 39362
$ws.Range("G51").Value2 = 19
$ws.Range("H51").Value2 = 1495

$ws.Range("B52").Value2 = 46867
$ws.Range("C52").Value2 = 290
$ws.Range("D52").Value2 = 33346
$ws.Range("E52").Value2 = 12571
$ws.Range("G52").Value2 = 5
$ws.Range("H52").Value2 = 950

$ws.Range("B53").Value2 = 44397
$ws.Range("C53").Value2 = 386
$ws.Range("D53").Value2 = 41209
$ws.Range("E53").Value2 = 3025

$ws.Range("B74").Value2 = 18494
$ws.Range("C74").Value2 = 141
$ws.Range("E74").Value2 = 5121

$ws.Range("B86").Value2 = 9684
$ws.Range("C86").Value2 = 46
$ws.Range("E86").Value2 = 571

$ws.Range("B89").Value2 = 8324
$ws.Range("C89").Value2 = 57
$ws.Range("D89").Value2 = 7549
$ws.Range("E89").Value2 = 726
$ws.Range("G89").Value2 = 1
$ws.Range("H89").Value2 = 49

$ws.Range("B99").Value2 = 6555
$ws.Range("C99").Value2 = 32
$ws.Range("D99").Value2 = 5570
$ws.Range("E99").Value2 = 828
$ws.Range("G99").Value2 = 0
$ws.Range("H99").Value2 = 157

$ws.Range("B100").Value2 = 6536
$ws.Range("C100").Value2 = 125
$ws.Range("D100").Value2 = 3379
$ws.Range("E100").Value2 = 2957
$ws.Range("G100").Value2 = 1
$ws.Range("H100").Value2 = 200

$ws.Range("B103").Value2 = 5541
$ws.Range("C103").Value2 = 90
$ws.Range("D103").Value2 = 710
$ws.Range("E103").Value2 = 4711
$ws.Range("G103").Value2 = 1
$ws.Range("H103").Value2 = 120

$ws.Range("B109").Value2 = 4674
$ws.Range("C109").Value2 = 16
$ws.Range("D109").Value2 = 2430
$ws.Range("E109").Value2 = 2098

$ws.Range("B113").Value2 = 3696
$ws.Range("C113").Value2 = 78
$ws.Range("D113").Value2 = 2521
$ws.Range("E113").Value2 = 1107
$ws.Range("G113").Value2 = 4
$ws.Range("H113").Value2 = 68

$ws.Range("B114").Value2 = 3664
$ws.Range("D114").Value2 = 1589
$ws.Range("E114").Value2 = 2017
$ws.Range("H114").Value2 = 58

$ws.Range("B116").Value2 = 3309
$ws.Range("C116").Value2 = 73
$ws.Range("D116").Value2 = 1634
$ws.Range("E116").Value2 = 1614
$ws.Range("G116").Value2 = 3
$ws.Range("H116").Value2 = 61

$ws.Range("B140").Value2 = 1364
$ws.Range("C140").Value2 = 11
$ws.Range("D140").Value2 = 1146
$ws.Range("E140").Value2 = 181

$ws.Range("B163").Value2 = 690
$ws.Range("C163").Value2 = 3
$ws.Range("E163").Value2 = 54

$ws.Range("B178").Value2 = 281
$ws.Range("C178").Value2 = 2
$ws.Range("E178").Value2 = 135

$ws.Range("B200").Value2 = 32
$ws.Range("C200").Value2 = 1
$ws.Range("E200").Value2 = 1

$ws.Range("D213").Value2 = 13
$ws.Range("H213").Value2 = 0

$ws.Range("D214").Value2 = 12
$ws.Range("H214").Value2 = 1

# --- Update "last updated" timestamp ---
$ws.Range("A1").Value2 = "Datos actualizados a 11 de Agosto de 2020 a las 01:24"
